# Rename header F1 from "BMI" to "Gain/Loss Amount"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F1").Value = "Gain/Loss Amount"

# --- New data rows (WeightID 1..7, dates 12/25/2017..12/31/2017 as serials) ---
$ws.Range("A2").Value = 1
$ws.Range("B2").Value = 43094
$ws.Range("C2").Value = 16.8
$ws.Range("D2").Value = 105
$ws.Range("E2").Value = 231
$ws.Range("F2").Value = 0

$ws.Range("A3").Value = 2
$ws.Range("B3").Value = 43095
$ws.Range("E3").Value = 0
$ws.Range("F3").Formula = "=E3-E2"

$ws.Range("A4").Value = 3
$ws.Range("B4").Value = 43096

$ws.Range("A5").Value = 4
$ws.Range("B5").Value = 43097

$ws.Range("A6").Value = 5
$ws.Range("B6").Value = 43098

$ws.Range("A7").Value = 6
$ws.Range("B7").Value = 43099

$ws.Range("A8").Value = 7
$ws.Range("B8").Value = 43100

# --- Number formats, applied in the same order Excel registers them (date,
# then the decimal weight format, then the WeightID integer format) so the
# generated numFmtId/cellXfs sequence lines up with the real workbook's. ---
$ws.Range("B2:B8").NumberFormat = "mm-dd-yy"
$ws.Range("C2:F3").NumberFormat = "0.0"
$ws.Range("A2:A8").NumberFormat = "0"

# --- Column widths for the newly populated columns ---
$ws.Columns("A").ColumnWidth = 8
$ws.Columns("B").ColumnWidth = 10
$ws.Columns("F").ColumnWidth = 15

# --- Selection after edits lands on A9, matching the saved workbook state ---
$ws.Range("A9").Select()
